$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.133.14"
$ws.Range("E2").Value = "  -2.18%  "
$ws.Range("D3").Value = "1.578.06"
$ws.Range("E3").Value = "  -1.45%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'209.19"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("E6").Value = "  -3.34%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  -0.86%  "
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("E10").Value = "  -0.87%  "
$ws.Range("D11").Value = "'0.0845"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D12").Value = "1.799.99"
$ws.Range("E12").Value = "  -1.44%  "
$ws.Range("D13").Value = "1.603.50"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("D15").Value = "'0.513"
$ws.Range("E15").Value = "  -1.88%  "
$ws.Range("D16").Value = "'64.43"
$ws.Range("E16").Value = "  -0.98%  "
$ws.Range("D17").Value = "26.145.69"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").Value = "'207.75"
$ws.Range("E21").Value = "  -1.23%  "
$ws.Range("E22").Value = "  -0.87%  "
$ws.Range("E23").Value = "  -2.39%  "
$ws.Range("E24").Value = "  -1.32%  "
$ws.Range("D25").Value = "'144.24"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("E27").Value = "  -1.56%  "
$ws.Range("E28").Value = "  -1.51%  "
$ws.Range("D29").Value = "'15.21"
$ws.Range("E29").Value = "  -1.04%  "
$ws.Range("D30").Value = "'0.0504"
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("E32").Value = "  -1.70%  "
$ws.Range("D33").Value = "'2.96"
$ws.Range("D34").Value = "1.277.14"
$ws.Range("E34").Value = "  -0.90%  "
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("D36").Value = "'0.609"
$ws.Range("E36").Value = "  +1.19%  "
$ws.Range("E37").Value = "  -1.03%  "
$ws.Range("D38").Value = "'1.13"
$ws.Range("E38").Value = "  -3.16%  "
$ws.Range("E39").Value = "  -2.25%  "
$ws.Range("D40").Value = "'0.815"
$ws.Range("E40").Value = "  -1.80%  "
$ws.Range("D41").Value = "'5.55"
$ws.Range("E41").Value = "  +2.85%  "
$ws.Range("E42").Value = "  -2.47%  "
$ws.Range("D43").Value = "'0.763"
$ws.Range("E43").Value = "  -2.52%  "
$ws.Range("D44").Value = "'62.43"
$ws.Range("E44").Value = "  -0.56%  "
$ws.Range("D45").Value = "1.713.54"
$ws.Range("E45").Value = "  -1.39%  "
$ws.Range("D46").Value = "'88.80"
$ws.Range("E46").Value = "  -1.91%  "
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("E48").Value = "  -1.63%  "
$ws.Range("D49").Value = "'0.100"
$ws.Range("E49").Value = "  -1.18%  "
$ws.Range("E50").Value = "  -2.17%  "
$ws.Range("E51").Value = "  -0.16%  "
